$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.920742638328592
$ws.Range("D2").Value = 9.667037446484882
$ws.Range("E2").Value = 13.75795227573954
$ws.Range("F2").Value = 30.9608575837635
$ws.Range("G2").Value = 33.36414497991077
$ws.Range("H2").Value = 14.37340777224141
$ws.Range("J2").Value = 9.732282694591145
$ws.Range("O2").Value = 22.74505173743368
$ws.Range("C3").Value = 4.74956410384177
$ws.Range("D3").Value = 9.673078951589259
$ws.Range("E3").Value = 13.72986351928126
$ws.Range("F3").Value = 30.60052972713553
$ws.Range("G3").Value = 32.53886467159906
$ws.Range("H3").Value = 14.3199452415204
$ws.Range("J3").Value = 9.734587955838133
$ws.Range("O3").Value = 22.51172477890332
$ws.Range("C4").Value = 4.642857908240282
$ws.Range("D4").Value = 9.678311898181548
$ws.Range("E4").Value = 13.71559754511369
$ws.Range("F4").Value = 30.38616264653964
$ws.Range("G4").Value = 32.03238208946689
$ws.Range("H4").Value = 14.2903671573651
$ws.Range("J4").Value = 9.737749318922226
$ws.Range("O4").Value = 22.37365562155597
$ws.Range("C5").Value = 4.599049569503739
$ws.Range("D5").Value = 9.680827631216973
$ws.Range("E5").Value = 13.71053700801866
$ws.Range("F5").Value = 30.30063607324186
$ws.Range("G5").Value = 31.82639236883861
$ws.Range("H5").Value = 14.27913879270443
$ws.Range("J5").Value = 9.739476756589823
$ws.Range("O5").Value = 22.31875788180995
$ws.Range("C6").Value = 4.591758020446794
$ws.Range("D6").Value = 9.681268518401264
$ws.Range("E6").Value = 13.70974227467619
$ws.Range("F6").Value = 30.2865478281831
$ws.Range("G6").Value = 31.79222273498842
$ws.Range("H6").Value = 14.27732437483258
$ws.Range("J6").Value = 9.739790120798043
$ws.Range("O6").Value = 22.30972639537169
$ws.Range("C7").Value = 4.642268300577557
$ws.Range("D7").Value = 9.67834427424569
$ws.Range("E7").Value = 13.71552624415198
$ws.Range("F7").Value = 30.38500166663745
$ws.Range("G7").Value = 32.02960193277704
$ws.Range("H7").Value = 14.2902123775646
$ws.Range("J7").Value = 9.737770837578884
$ws.Range("O7").Value = 22.37290964086506
$ws.Range("C8").Value = 4.862097684253725
$ws.Range("D8").Value = 9.668804505872895
$ws.Range("E8").Value = 13.74765091334335
$ws.Range("F8").Value = 30.8352511758388
$ws.Range("G8").Value = 33.07974150062982
$ws.Range("H8").Value = 14.35430495193123
$ws.Range("J8").Value = 9.732715306555376
$ws.Range("O8").Value = 22.66355957609658
$ws.Range("C9").Value = 5.277361014224274
$ws.Range("D9").Value = 9.662173792963259
$ws.Range("E9").Value = 13.83411694412124
$ws.Range("F9").Value = 31.76813656619898
$ws.Range("G9").Value = 35.12644522169583
$ws.Range("H9").Value = 14.50535820220195
$ws.Range("J9").Value = 9.736643243213058
$ws.Range("O9").Value = 23.2719534676994
$ws.Range("C10").Value = 5.569097360695009
$ws.Range("D10").Value = 9.664645213267436
$ws.Range("E10").Value = 13.9116687983004
$ws.Range("F10").Value = 32.4776198614901
$ws.Range("G10").Value = 36.60406146725428
$ws.Range("H10").Value = 14.63119705848083
$ws.Range("J10").Value = 9.747943953925208
$ws.Range("O10").Value = 23.73848429676602
$ws.Range("C11").Value = 5.69825738323224
$ws.Range("D11").Value = 9.66735770262812
$ws.Range("E11").Value = 13.94992080461271
$ws.Range("F11").Value = 32.80423388496732
$ws.Range("G11").Value = 37.26706215651704
$ws.Range("H11").Value = 14.69152319190209
$ws.Range("J11").Value = 9.754903559497089
$ws.Range("O11").Value = 23.95410895876937
$ws.Range("C12").Value = 5.746609374124008
$ws.Range("D12").Value = 9.668612479726315
$ws.Range("E12").Value = 13.96482610652025
$ws.Range("F12").Value = 32.92835452901541
$ws.Range("G12").Value = 37.51653606913305
$ws.Range("H12").Value = 14.71479642540413
$ws.Range("J12").Value = 9.757799456284602
$ws.Range("O12").Value = 24.03617626535203
$ws.Range("C13").Value = 5.73622142367559
$ws.Range("D13").Value = 9.668332131241547
$ws.Range("E13").Value = 13.96159742434258
$ws.Range("F13").Value = 32.9016052534456
$ws.Range("G13").Value = 37.46288220902905
$ws.Range("H13").Value = 14.70976526868477
$ws.Range("J13").Value = 9.757164209302688
$ws.Range("O13").Value = 24.0184843198436
$ws.Range("C14").Value = 5.702246767003851
$ws.Range("D14").Value = 9.66745637695263
$ws.Range("E14").Value = 13.95113869827606
$ws.Range("F14").Value = 32.81443723296751
$ws.Range("G14").Value = 37.28761987037166
$ws.Range("H14").Value = 14.69342937958463
$ws.Range("J14").Value = 9.755136593949459
$ws.Range("O14").Value = 23.96085279118099
$ws.Range("C15").Value = 5.681362318579184
$ws.Range("D15").Value = 9.666949569670352
$ws.Range("E15").Value = 13.94478690389392
$ws.Range("F15").Value = 32.76109802074144
$ws.Range("G15").Value = 37.18005186979804
$ws.Range("H15").Value = 14.68347861846556
$ws.Range("J15").Value = 9.75392850268393
$ws.Range("O15").Value = 23.92560366731032
$ws.Range("C16").Value = 5.560580797131274
$ws.Range("D16").Value = 9.664499841887627
$ws.Range("E16").Value = 13.90922814636433
$ws.Range("F16").Value = 32.456342647335
$ws.Range("G16").Value = 36.5605245828022
$ws.Range("H16").Value = 14.6273154247838
$ws.Range("J16").Value = 9.747525637272311
$ws.Range("O16").Value = 23.72445479040498
$ws.Range("C17").Value = 5.485539969987388
$ws.Range("D17").Value = 9.66340328699247
$ws.Range("E17").Value = 13.88817041975872
$ws.Range("F17").Value = 32.27029187686014
$ws.Range("G17").Value = 36.17791568812815
$ws.Range("H17").Value = 14.59364058654361
$ws.Range("J17").Value = 9.744062813389009
$ws.Range("O17").Value = 23.60187428054169
$ws.Range("C18").Value = 5.442046567105601
$ws.Range("D18").Value = 9.662922122700735
$ws.Range("E18").Value = 13.87633891779286
$ws.Range("F18").Value = 32.16365129537095
$ws.Range("G18").Value = 35.95700053572275
$ws.Range("H18").Value = 14.57456264014794
$ws.Range("J18").Value = 9.74224236600921
$ws.Range("O18").Value = 23.53169342163405
$ws.Range("C19").Value = 5.427264976841309
$ws.Range("D19").Value = 9.6627849147264
$ws.Range("E19").Value = 13.8723813411688
$ws.Range("F19").Value = 32.12761190950916
$ws.Range("G19").Value = 35.88206529527471
$ws.Range("H19").Value = 14.56815355859894
$ws.Range("J19").Value = 9.741655441399795
$ws.Range("O19").Value = 23.50798931408281
$ws.Range("C20").Value = 5.493562905250963
$ws.Range("D20").Value = 9.663504544996401
$ws.Range("E20").Value = 13.8903830880996
$ws.Range("F20").Value = 32.29005976310948
$ws.Range("G20").Value = 36.21873491593509
$ws.Range("H20").Value = 14.59719531999442
$ws.Range("J20").Value = 9.744413717773387
$ws.Range("O20").Value = 23.61489016423073
$ws.Range("C21").Value = 5.712241443374235
$ws.Range("D21").Value = 9.667707436493219
$ws.Range("E21").Value = 13.95419934077838
$ws.Range("F21").Value = 32.84002960602479
$ws.Range("G21").Value = 37.33914387723872
$ws.Range("H21").Value = 14.69821609815423
$ws.Range("J21").Value = 9.755725095049417
$ws.Range("O21").Value = 23.97776989354938
$ws.Range("C22").Value = 5.85188934449856
$ws.Range("D22").Value = 9.671780550157157
$ws.Range("E22").Value = 13.99835171587853
$ws.Range("F22").Value = 33.20197221035583
$ws.Range("G22").Value = 38.06201678356244
$ws.Range("H22").Value = 14.7667327204538
$ws.Range("J22").Value = 9.764635163502575
$ws.Range("O22").Value = 24.21731728487537
$ws.Range("C23").Value = 5.777670146970701
$ws.Range("D23").Value = 9.669485582652015
$ws.Range("E23").Value = 13.9745656859514
$ws.Range("F23").Value = 33.00860570262192
$ws.Range("G23").Value = 37.67714822793401
$ws.Range("H23").Value = 14.72994077576603
$ws.Range("J23").Value = 9.759741260557602
$ws.Range("O23").Value = 24.08927243296216
$ws.Range("C24").Value = 5.489936827226254
$ws.Range("D24").Value = 9.663458301277773
$ws.Range("E24").Value = 13.88938188403092
$ws.Range("F24").Value = 32.28112168827865
$ws.Range("G24").Value = 36.20028347563822
$ws.Range("H24").Value = 14.59558734499336
$ws.Range("J24").Value = 9.744254543125525
$ws.Range("O24").Value = 23.60900476706775
$ws.Range("C25").Value = 5.167123255842992
$ws.Range("D25").Value = 9.662676225994488
$ws.Range("E25").Value = 13.80823951132447
$ws.Range("F25").Value = 31.51106847167174
$ws.Range("G25").Value = 34.57607532276457
$ws.Range("H25").Value = 14.46183719371571
$ws.Range("J25").Value = 9.734100479248381
$ws.Range("O25").Value = 23.10364695539334
